# Re-curate the dimension metadata for the "municipio-nombre" column (M).
# It now follows the same pattern already used by "provincia-nombre" (Q)
# and "comarca-nombre" (U): an sdmx-dimension:refArea dimension, typed as
# "dim", whose concept is identified by its own URI column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "sdmx-dimension:refArea"
$ws.Range("M3").Value = "dim"
$ws.Range("M4").Value = "URI-Municipio"
